# crassets.xlsx: add new symbol TIA22861-USD (and the other two new rows
# that came with it: OSMO, AGIX) to the portfolio sheet, and correct the
# B6 (ATOM "Anzahl") value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row tweak: ATOM amount changed from 1488 to 1393.
$ws.Range("B6").Value = 1393

# New row 28: OSMO
$ws.Range("A28").Value = "OSMO"
$ws.Range("B28").Value = 621
$ws.Range("C28").Value = 1.7
$ws.Range("D28").Value = "DeFi/So"

# New row 29: TIA
$ws.Range("A29").Value = "TIA"
$ws.Range("B29").Value = 27.5
$ws.Range("C29").Value = 19
$ws.Range("D29").Value = "Utility"

# New row 30: AGIX
$ws.Range("A30").Value = "AGIX"
$ws.Range("B30").Value = 1232
$ws.Range("C30").Value = 0.435
$ws.Range("D30").Value = "KI"

# Match the author's final selection in the saved workbook.
$ws.Range("H13").Select()
